# Update odds values in row 2 and row 3 per the FlashScore data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.5
$ws.Range("H2").Value = 2.95
$ws.Range("J2").Value = 3.15
$ws.Range("K2").Value = 1.98
$ws.Range("L2").Value = 3.25
$ws.Range("N2").Value = 6.5
$ws.Range("P2").Value = 2.94
$ws.Range("Q2").Value = 1.88
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 2.42
$ws.Range("U2").Value = 1.73
$ws.Range("V2").Value = 2.06
$ws.Range("W2").Value = 6.8
$ws.Range("X2").Value = 10.5
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 17.5
$ws.Range("AB2").Value = 23
$ws.Range("AC2").Value = 8.5
$ws.Range("AH2").Value = 7.8
$ws.Range("AI2").Value = 12.5
$ws.Range("AL2").Value = 17.5
$ws.Range("AM2").Value = 22
$ws.Range("AN2").Value = 4.45
$ws.Range("AO2").Value = 14.5
$ws.Range("AP2").Value = 22
$ws.Range("AQ2").Value = 65
$ws.Range("AR2").Value = 100
$ws.Range("AS2").Value = 300
$ws.Range("AT2").Value = 2.42
$ws.Range("AU2").Value = 6.6
$ws.Range("AV2").Value = 55
$ws.Range("AW2").Value = 4.65
$ws.Range("AX2").Value = 14.5
$ws.Range("AY2").Value = 20
$ws.Range("AZ2").Value = 65
$ws.Range("BA2").Value = 90
$ws.Range("BB2").Value = 250
$ws.Range("G3").Value = 2.75
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 2.45
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 3.2
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("W3").Value = 9.5
$ws.Range("X3").Value = 13
$ws.Range("Z3").Value = 29
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 29
$ws.Range("AE3").Value = 13
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 12
$ws.Range("AJ3").Value = 10
$ws.Range("AK3").Value = 23
$ws.Range("AL3").Value = 21
$ws.Range("AN3").Value = 4.75
$ws.Range("AO3").Value = 15
$ws.Range("AP3").Value = 23
$ws.Range("AR3").Value = 67
$ws.Range("AW3").Value = 4.5
$ws.Range("BA3").Value = 67
